$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.67'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '24.93'
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.185'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05733'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.479'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.059'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8099'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8416'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1338'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06992'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02814'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09373'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001510'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0005958'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '14OneONE'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006232'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.499'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.122'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3195'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03208'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1301'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.745'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04667'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1329'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001236'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004254'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009693'
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001501'
$ws.Range("D28").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03630'
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = 'BKEXToken'

$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1050'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '40BKEXTokenBKK'

$ws.Range("B42").Value = 'CEJI'

$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002738'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '41CEJICEJIBestin24h'

$ws.Range("B43").Value = 'KickToken'

$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003235'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '42KickTokenKICK'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007312'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005285'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1899'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("D50").Style = "Normal"
